$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $value) {
    $c = $ws.Range($ref)
    $c.Value = "'" + $value
    $c.Style = "Normal"
}

Set-TextValue "D2" "43.650.78"
Set-TextValue "E2" "  +1.21%  "
Set-TextValue "D3" "2.417.16"
Set-TextValue "E3" "  +2.37%  "
Set-TextValue "D4" "1.00"
Set-TextValue "E4" "  +0.04%  "
Set-TextValue "D5" "306.78"
Set-TextValue "E5" "  +1.23%  "
Set-TextValue "E6" "  +1.11%  "
Set-TextValue "E7" "  +0.44%  "
Set-TextValue "E8" "  +0.04%  "
Set-TextValue "E9" "  -1.30%  "
Set-TextValue "D10" "35.21"
Set-TextValue "E10" "  +2.97%  "
Set-TextValue "E11" "  +1.19%  "
Set-TextValue "E12" "  +2.61%  "
Set-TextValue "D13" "18.54"
Set-TextValue "E13" "  -0.66%  "
Set-TextValue "D14" "6.90"
Set-TextValue "E14" "  +2.29%  "
Set-TextValue "D15" "2.783.06"
Set-TextValue "E15" "  +2.22%  "
Set-TextValue "D16" "2.443.66"
Set-TextValue "E16" "  +1.38%  "
Set-TextValue "D17" "0.826"
Set-TextValue "E17" "  +3.57%  "
Set-TextValue "D18" "43.629.37"
Set-TextValue "E18" "  +1.21%  "
Set-TextValue "E19" "  +2.40%  "
Set-TextValue "E20" "  -0.98%  "
Set-TextValue "D21" "0.0₃0902"
Set-TextValue "E21" "  +1.43%  "
Set-TextValue "D22" "68.34"
Set-TextValue "E22" "  +0.24%  "
Set-TextValue "D23" "237.94"
Set-TextValue "E23" "  +0.91%  "
Set-TextValue "E24" "  +0.69%  "
Set-TextValue "E25" "  +0.84%  "
Set-TextValue "E26" "  +0.07%  "
Set-TextValue "D27" "24.98"
Set-TextValue "E27" "  +1.69%  "
Set-TextValue "E28" "  -0.60%  "
Set-TextValue "D29" "9.44"
Set-TextValue "E29" "  +3.42%  "
Set-TextValue "D30" "32.48"
Set-TextValue "E30" "  +3.49%  "
Set-TextValue "D31" "0.118"
Set-TextValue "E31" "  +16.65%  "
Set-TextValue "E32" "  +7.10%  "
Set-TextValue "D33" "5.14"
Set-TextValue "E33" "  +1.50%  "
Set-TextValue "E35" "  +3.39%  "
Set-TextValue "E36" "  +3.08%  "
Set-TextValue "D37" "130.66"
Set-TextValue "E37" "  +24.61%  "
Set-TextValue "D38" "2.92"
Set-TextValue "E38" "  +5.51%  "
Set-TextValue "D39" "4.41"
Set-TextValue "E39" "  +0.31%  "
Set-TextValue "E40" "  -1.07%  "
Set-TextValue "E41" "  -0.12%  "
Set-TextValue "D42" "21.38"
Set-TextValue "E42" "  -4.64%  "
Set-TextValue "D43" "1.948.71"
Set-TextValue "E43" "  +0.20%  "
Set-TextValue "E44" "  +1.31%  "
Set-TextValue "E45" "  +2.12%  "
Set-TextValue "E46" "  +3.19%  "
Set-TextValue "D47" "9.33"
Set-TextValue "E47" "  -0.96%  "
Set-TextValue "D48" "2.636.22"
Set-TextValue "E48" "  +2.05%  "
Set-TextValue "D49" "1.57"
Set-TextValue "E49" "  +3.95%  "
Set-TextValue "D50" "52.80"
Set-TextValue "E50" "  -0.44%  "
Set-TextValue "D51" "72.37"
Set-TextValue "E51" "  +0.03%  "
